# Append a new row (21) of station-load data to the bottom of the sheet,
# matching the "actual energy file" refresh described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 21

# Column A holds the date as plain text (e.g. "2025-01-18"), same as every
# other row above it. Force a text number format first so Excel doesn't
# auto-convert the recognisable date-like string into a date serial value,
# then restore the Normal style so the cell doesn't pick up a "number
# stored as text" quote-prefix style (matching the plain formatting used
# by every other row in this column).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-01-18"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 20
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 15333
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 15333
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 15334.5333
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 18).Value = 0
$ws.Cells.Item($row, 19).Value = 0
$ws.Cells.Item($row, 20).Value = 15334.5333
